$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# The first available name in the pool gets "used" - grab it before removing it.
$usedId = $namesSheet.Range("A1").Value2

# Remove the consumed name from the pool, shifting everything up.
$namesSheet.Rows.Item(1).Delete()

# Append a new row to the "used" log for this name.
$nextRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row + 1
$usedSheet.Cells.Item($nextRow, 1).Value = $usedId
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月18日 08_44_52.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-18 08:51:37"
